$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: Original Product category
# Column B: Complementary Product category
# Column C: frequency

$ws.Range("A2").Value = "Kitchen & Dining"
$ws.Range("B2").Value = "Kitchen & Dining"
$ws.Range("C2").Value = 783

$ws.Range("A3").Value = "Fashion & Accessories"
$ws.Range("B3").Value = "Fashion & Accessories"
$ws.Range("C3").Value = 562

$ws.Range("A4").Value = "Storage & Organization"
$ws.Range("B4").Value = "Storage & Organization"
$ws.Range("C4").Value = 316

$ws.Range("A5").Value = "Kids & Toys"
$ws.Range("B5").Value = "Storage & Organization"
$ws.Range("C5").Value = 146

$ws.Range("A6").Value = "Fashion & Accessories"
$ws.Range("B6").Value = "Storage & Organization"
$ws.Range("C6").Value = 103

$ws.Range("A7").Value = "Storage & Organization"
$ws.Range("B7").Value = "Fashion & Accessories"
$ws.Range("C7").Value = 91

$ws.Range("A8").Value = "Storage & Organization"
$ws.Range("B8").Value = "Kids & Toys"
$ws.Range("C8").Value = 80

$ws.Range("A9").Value = "Fashion & Accessories"
$ws.Range("B9").Value = "Kitchen & Dining"
$ws.Range("C9").Value = 65

$ws.Range("A10").Value = "Kids & Toys"
$ws.Range("B10").Value = "Kids & Toys"
$ws.Range("C10").Value = 63

$ws.Range("A11").Value = "Home Decor"
$ws.Range("B11").Value = "Home Decor"
$ws.Range("C11").Value = 62
